$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 90-113 with refreshed rolling-window data (one new day inserted on 2021-02-08) ---
$ws.Cells.Item(90, 1).Value = 44232
$ws.Cells.Item(90, 2).Value = 78
$ws.Cells.Item(90, 3).Value = 470
$ws.Cells.Item(90, 4).Value = 248.6601450693868

$ws.Cells.Item(91, 1).Value = 44233
$ws.Cells.Item(91, 2).Value = 78
$ws.Cells.Item(91, 3).Value = 455
$ws.Cells.Item(91, 4).Value = 240.7241829927042

$ws.Cells.Item(92, 1).Value = 44234
$ws.Cells.Item(92, 2).Value = 77
$ws.Cells.Item(92, 3).Value = 453
$ws.Cells.Item(92, 4).Value = 239.6660547158132

$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 117
$ws.Cells.Item(93, 3).Value = 480
$ws.Cells.Item(93, 4).Value = 253.9507864538418

$ws.Cells.Item(94, 1).Value = 44236
$ws.Cells.Item(94, 2).Value = 31
$ws.Cells.Item(94, 3).Value = 474
$ws.Cells.Item(94, 4).Value = 250.7764016231688

$ws.Cells.Item(95, 1).Value = 44237
$ws.Cells.Item(95, 2).Value = 14
$ws.Cells.Item(95, 3).Value = 488
$ws.Cells.Item(95, 4).Value = 258.1832995614058

$ws.Cells.Item(96, 1).Value = 44238
$ws.Cells.Item(96, 2).Value = 85
$ws.Cells.Item(96, 3).Value = 480
$ws.Cells.Item(96, 4).Value = 253.9507864538418

$ws.Cells.Item(97, 1).Value = 44239
$ws.Cells.Item(97, 2).Value = 72
$ws.Cells.Item(97, 3).Value = 426
$ws.Cells.Item(97, 4).Value = 225.3813229777846

$ws.Cells.Item(98, 1).Value = 44240
$ws.Cells.Item(98, 2).Value = 92
$ws.Cells.Item(98, 3).Value = 436
$ws.Cells.Item(98, 4).Value = 230.6719643622396

$ws.Cells.Item(99, 1).Value = 44241
$ws.Cells.Item(99, 2).Value = 69
$ws.Cells.Item(99, 3).Value = 443
$ws.Cells.Item(99, 4).Value = 234.3754133313581

$ws.Cells.Item(100, 1).Value = 44242
$ws.Cells.Item(100, 2).Value = 63
$ws.Cells.Item(100, 3).Value = 413
$ws.Cells.Item(100, 4).Value = 218.5034891779931

$ws.Cells.Item(101, 1).Value = 44243
$ws.Cells.Item(101, 2).Value = 41
$ws.Cells.Item(101, 3).Value = 399
$ws.Cells.Item(101, 4).Value = 211.096591239756

$ws.Cells.Item(102, 1).Value = 44244
$ws.Cells.Item(102, 2).Value = 21
$ws.Cells.Item(102, 3).Value = 373
$ws.Cells.Item(102, 4).Value = 197.3409236401729

$ws.Cells.Item(103, 1).Value = 44245
$ws.Cells.Item(103, 2).Value = 55
$ws.Cells.Item(103, 3).Value = 397
$ws.Cells.Item(103, 4).Value = 210.038462962865

$ws.Cells.Item(104, 1).Value = 44246
$ws.Cells.Item(104, 2).Value = 58
$ws.Cells.Item(104, 3).Value = 405
$ws.Cells.Item(104, 4).Value = 214.270976070429

$ws.Cells.Item(105, 1).Value = 44247
$ws.Cells.Item(105, 2).Value = 66
$ws.Cells.Item(105, 3).Value = 407
$ws.Cells.Item(105, 4).Value = 215.32910434732

$ws.Cells.Item(106, 1).Value = 44248
$ws.Cells.Item(106, 2).Value = 93
$ws.Cells.Item(106, 3).Value = 398
$ws.Cells.Item(106, 4).Value = 210.5675271013105

$ws.Cells.Item(107, 1).Value = 44249
$ws.Cells.Item(107, 2).Value = 71
$ws.Cells.Item(107, 3).Value = 438
$ws.Cells.Item(107, 4).Value = 231.7300926391306

$ws.Cells.Item(108, 1).Value = 44250
$ws.Cells.Item(108, 2).Value = 43
$ws.Cells.Item(108, 3).Value = 504
$ws.Cells.Item(108, 4).Value = 266.6483257765339

$ws.Cells.Item(109, 1).Value = 44251
$ws.Cells.Item(109, 2).Value = 12
$ws.Cells.Item(109, 3).Value = 550
$ws.Cells.Item(109, 4).Value = 290.9852761450271

$ws.Cells.Item(110, 1).Value = 44252
$ws.Cells.Item(110, 2).Value = 95
$ws.Cells.Item(110, 3).Value = 567
$ws.Cells.Item(110, 4).Value = 299.9793664986006

$ws.Cells.Item(111, 1).Value = 44253
$ws.Cells.Item(111, 2).Value = 124
$ws.Cells.Item(111, 3).Value = 615
$ws.Cells.Item(111, 4).Value = 325.3744451439848

$ws.Cells.Item(112, 1).Value = 44254
$ws.Cells.Item(112, 2).Value = 112
$ws.Cells.Item(112, 3).Value = 623
$ws.Cells.Item(112, 4).Value = 329.6069582515489

$ws.Cells.Item(113, 1).Value = 44255
$ws.Cells.Item(113, 2).Value = 110

# --- Append brand new rows 114-115, copying formatting from the last existing row (113) ---
$ws.Range("A113:D113").Copy()
$ws.Range("A114:D115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(114, 1).Value = 44256
$ws.Cells.Item(114, 2).Value = 119
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 51

$ws.Range("A1").Select() | Out-Null
